# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.789.68"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").Value = "2.092.95"
$ws.Range("E3").Value = "  +2.54%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'228.27"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("D7").Value = "'60.69"
$ws.Range("E7").Value = "  +1.81%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").Value = "'0.0838"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "2.402.31"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").Value = "'14.97"
$ws.Range("E13").Value = "  +3.75%  "

$ws.Range("D14").Value = "'21.95"
$ws.Range("E14").Value = "  +4.35%  "

$ws.Range("D15").Value = "'0.797"

$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "2.087.41"
$ws.Range("E17").Value = "  +2.27%  "

$ws.Range("D18").Value = "38.749.47"
$ws.Range("E18").Value = "  +2.66%  "

$ws.Range("D19").Value = "'71.74"

$ws.Range("D20").Value = "'6.02"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("D22").Value = "'226.43"

$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "'170.47"
$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").Value = "'9.43"
$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").Value = "'0.137"
$ws.Range("E28").Value = "  +7.01%  "

$ws.Range("D29").Value = "'1.45"
$ws.Range("E29").Value = "  +12.28%  "

$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").Value = "'2.35"
$ws.Range("E32").Value = "  +4.63%  "

$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  +2.50%  "

$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  +4.95%  "

$ws.Range("D35").Value = "'0.0613"
$ws.Range("E35").Value = "  +2.01%  "

$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.43"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'18.20"
$ws.Range("E40").Value = "  +1.21%  "

$ws.Range("D41").Value = "1.538.87"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").Value = "'101.06"
$ws.Range("E42").Value = "  +3.69%  "

$ws.Range("D43").Value = "'0.0223"
$ws.Range("E43").Value = "  +3.51%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0926"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.82"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").Value = "'7.65"
$ws.Range("E46").Value = "  +8.13%  "

$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").Value = "2.287.80"
$ws.Range("E51").Value = "  +2.48%  "
